$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 959
$ws.Cells.Item(4, 9).Value = 901.25
$ws.Cells.Item(4, 10).Value = 1190
$ws.Cells.Item(4, 11).Value = 901.25
$ws.Cells.Item(4, 12).Value = 1190
$ws.Cells.Item(4, 13).Value = -787.25
$ws.Cells.Item(4, 14).Value = -1418
$ws.Cells.Item(9, 8).Value = 2502
$ws.Cells.Item(9, 9).Value = 471.6
$ws.Cells.Item(9, 10).Value = 5040
$ws.Cells.Item(9, 11).Value = 471.6
$ws.Cells.Item(9, 12).Value = 5040
$ws.Cells.Item(9, 13).Value = -302.6
$ws.Cells.Item(9, 14).Value = -5378
$ws.Cells.Item(10, 8).Value = 0
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 14).Value = $null
$ws.Cells.Item(21, 8).Value = 10017
$ws.Cells.Item(21, 9).Value = 10017
$ws.Cells.Item(21, 11).Value = 10017
$ws.Cells.Item(21, 13).Value = -9549
$ws.Cells.Item(23, 8).Value = 10017
$ws.Cells.Item(23, 9).Value = 10017
$ws.Cells.Item(23, 11).Value = 10017
$ws.Cells.Item(23, 13).Value = -9783
$ws.Cells.Item(29, 8).Value = 341.75
$ws.Cells.Item(29, 9).Value = 341.75
$ws.Cells.Item(29, 11).Value = 1025.25
$ws.Cells.Item(29, 13).Value = -744.25
$ws.Cells.Item(33, 8).Value = 469
$ws.Cells.Item(33, 9).Value = 524.25
$ws.Cells.Item(33, 11).Value = 524.25
$ws.Cells.Item(33, 13).Value = -295.25
$ws.Cells.Item(39, 8).Value = 50.857143
$ws.Cells.Item(39, 9).Value = 36.833332
$ws.Cells.Item(39, 11).Value = 110.499996
$ws.Cells.Item(39, 13).Value = 185.500004
$ws.Cells.Item(40, 8).Value = 5801.4287
$ws.Cells.Item(40, 9).Value = 3499
$ws.Cells.Item(40, 10).Value = 6185.1665
$ws.Cells.Item(40, 11).Value = 3499
$ws.Cells.Item(40, 12).Value = 6185.1665
$ws.Cells.Item(40, 13).Value = -3324
$ws.Cells.Item(40, 14).Value = -6535.1665
$ws.Cells.Item(55, 8).Value = 482.64706
$ws.Cells.Item(55, 9).Value = 434.3846
$ws.Cells.Item(55, 11).Value = 434.3846
$ws.Cells.Item(55, 13).Value = -220.3846
$ws.Cells.Item(58, 8).Value = 565.44446
$ws.Cells.Item(58, 10).Value = 983
$ws.Cells.Item(58, 12).Value = 2949
$ws.Cells.Item(58, 14).Value = -3249
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 14).Value = $null
$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 14).Value = $null
$ws.Cells.Item(81, 8).Value = 50000
$ws.Cells.Item(81, 10).Value = 50000
$ws.Cells.Item(81, 12).Value = 50000
$ws.Cells.Item(81, 14).Value = -51996
$ws.Cells.Item(84, 8).Value = 50000
$ws.Cells.Item(84, 10).Value = 50000
$ws.Cells.Item(84, 12).Value = 150000
$ws.Cells.Item(84, 14).Value = -159984
$ws.Cells.Item(87, 8).Value = 84494.5
$ws.Cells.Item(87, 10).Value = 84494.5
$ws.Cells.Item(87, 12).Value = 84494.5
$ws.Cells.Item(87, 14).Value = -86990.5
$ws.Cells.Item(90, 8).Value = 84494.5
$ws.Cells.Item(90, 10).Value = 84494.5
$ws.Cells.Item(90, 12).Value = 253483.5
$ws.Cells.Item(90, 14).Value = -265963.5
$ws.Cells.Item(98, 8).Value = 2908.4443
$ws.Cells.Item(98, 9).Value = 2908.4443
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 2908.4443
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = -1410.4443
$ws.Cells.Item(98, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 2908.4443
$ws.Cells.Item(122, 9).Value = 2908.4443
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 8725.332900000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -6275.332900000001
$ws.Cells.Item(122, 14).Value = $null
$ws.Cells.Item(128, 8).Value = 60000
$ws.Cells.Item(128, 10).Value = 60000
$ws.Cells.Item(128, 12).Value = 60000
$ws.Cells.Item(128, 14).Value = -69960
$ws.Cells.Item(131, 8).Value = 14949.5
$ws.Cells.Item(131, 9).Value = 14899
$ws.Cells.Item(131, 11).Value = 44697
$ws.Cells.Item(131, 13).Value = -39657
$ws.Cells.Item(132, 8).Value = 5931.0557
$ws.Cells.Item(132, 9).Value = 6320.7334
$ws.Cells.Item(132, 10).Value = 3982.6667
$ws.Cells.Item(132, 11).Value = 18962.2002
$ws.Cells.Item(132, 12).Value = 11948.0001
$ws.Cells.Item(132, 13).Value = -16432.2002
$ws.Cells.Item(132, 14).Value = -17008.0001
$ws.Cells.Item(137, 8).Value = 3128.75
$ws.Cells.Item(137, 9).Value = 2344.4
$ws.Cells.Item(137, 10).Value = 5089.625
$ws.Cells.Item(137, 11).Value = 7033.200000000001
$ws.Cells.Item(137, 12).Value = 15268.875
$ws.Cells.Item(137, 13).Value = -4483.200000000001
$ws.Cells.Item(137, 14).Value = -20368.875
$ws.Cells.Item(138, 8).Value = 3582.5557
$ws.Cells.Item(138, 10).Value = 4873.25
$ws.Cells.Item(138, 12).Value = 14619.75
$ws.Cells.Item(138, 14).Value = -24899.75
$ws.Cells.Item(141, 8).Value = 5228.278
$ws.Cells.Item(141, 9).Value = 4471.4116
$ws.Cells.Item(141, 11).Value = 13414.2348
$ws.Cells.Item(141, 13).Value = -8234.234800000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 357.81818
$ws.Cells.Item(5, 9).Value = 393.5
$ws.Cells.Item(5, 10).Value = 262.66666
$ws.Cells.Item(5, 11).Value = 393.5
$ws.Cells.Item(5, 12).Value = 262.66666
$ws.Cells.Item(5, 13).Value = -281.5
$ws.Cells.Item(5, 14).Value = -486.66666
$ws.Cells.Item(6, 8).Value = 666850
$ws.Cells.Item(6, 9).Value = 666850
$ws.Cells.Item(6, 11).Value = 666850
$ws.Cells.Item(6, 13).Value = -666677
$ws.Cells.Item(16, 8).Value = 1200
$ws.Cells.Item(16, 9).Value = 1200
$ws.Cells.Item(16, 11).Value = 1200
$ws.Cells.Item(16, 13).Value = -913
$ws.Cells.Item(31, 8).Value = 8823.333000000001
$ws.Cells.Item(31, 9).Value = 8823.333000000001
$ws.Cells.Item(31, 11).Value = 8823.333000000001
$ws.Cells.Item(31, 13).Value = -8529.333000000001
$ws.Cells.Item(50, 8).Value = 2490
$ws.Cells.Item(50, 9).Value = 2500
$ws.Cells.Item(50, 10).Value = 2485
$ws.Cells.Item(50, 11).Value = 2500
$ws.Cells.Item(50, 12).Value = 2485
$ws.Cells.Item(50, 13).Value = -1786
$ws.Cells.Item(50, 14).Value = -3913
$ws.Cells.Item(88, 8).Value = 2144
$ws.Cells.Item(88, 9).Value = 1247.5
$ws.Cells.Item(88, 10).Value = 2442.8333
$ws.Cells.Item(88, 11).Value = 1247.5
$ws.Cells.Item(88, 12).Value = 2442.8333
$ws.Cells.Item(88, 13).Value = -841.5
$ws.Cells.Item(88, 14).Value = -3254.8333
$ws.Cells.Item(91, 8).Value = 2144
$ws.Cells.Item(91, 9).Value = 1247.5
$ws.Cells.Item(91, 10).Value = 2442.8333
$ws.Cells.Item(91, 11).Value = 1247.5
$ws.Cells.Item(91, 12).Value = 2442.8333
$ws.Cells.Item(91, 13).Value = 156.5
$ws.Cells.Item(91, 14).Value = -5250.8333
$ws.Cells.Item(132, 8).Value = 3338.4333
$ws.Cells.Item(132, 9).Value = 3430.6785
$ws.Cells.Item(132, 11).Value = 10292.0355
$ws.Cells.Item(132, 13).Value = -7762.0355

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 357.81818
$ws.Cells.Item(4, 9).Value = 393.5
$ws.Cells.Item(4, 10).Value = 262.66666
$ws.Cells.Item(4, 11).Value = 393.5
$ws.Cells.Item(4, 12).Value = 262.66666
$ws.Cells.Item(4, 13).Value = -278.5
$ws.Cells.Item(4, 14).Value = -492.66666
$ws.Cells.Item(11, 8).Value = 926.3
$ws.Cells.Item(11, 9).Value = 329
$ws.Cells.Item(11, 10).Value = 1822.25
$ws.Cells.Item(11, 11).Value = 329
$ws.Cells.Item(11, 12).Value = 1822.25
$ws.Cells.Item(11, 13).Value = -189
$ws.Cells.Item(11, 14).Value = -2102.25
$ws.Cells.Item(58, 8).Value = 62992.25
$ws.Cells.Item(58, 10).Value = 62992.25
$ws.Cells.Item(58, 12).Value = 62992.25
$ws.Cells.Item(58, 14).Value = -63580.25
$ws.Cells.Item(74, 8).Value = 80620
$ws.Cells.Item(74, 10).Value = 80620
$ws.Cells.Item(74, 12).Value = 80620
$ws.Cells.Item(74, 14).Value = -82492
$ws.Cells.Item(77, 8).Value = 80620
$ws.Cells.Item(77, 10).Value = 80620
$ws.Cells.Item(77, 12).Value = 241860
$ws.Cells.Item(77, 14).Value = -251220
$ws.Cells.Item(80, 8).Value = 1232.6154
$ws.Cells.Item(80, 10).Value = 2029.8572
$ws.Cells.Item(80, 12).Value = 2029.8572
$ws.Cells.Item(80, 14).Value = -4025.8572
$ws.Cells.Item(83, 8).Value = 1232.6154
$ws.Cells.Item(83, 10).Value = 2029.8572
$ws.Cells.Item(83, 12).Value = 10149.286
$ws.Cells.Item(83, 14).Value = -20133.286
$ws.Cells.Item(99, 8).Value = 6169.8
$ws.Cells.Item(99, 9).Value = 4712.25
$ws.Cells.Item(99, 11).Value = 4712.25
$ws.Cells.Item(99, 13).Value = -3214.25
$ws.Cells.Item(134, 8).Value = 4466.75
$ws.Cells.Item(134, 9).Value = 3624.2
$ws.Cells.Item(134, 11).Value = 10872.6
$ws.Cells.Item(134, 13).Value = -8337.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1703.3334
$ws.Cells.Item(16, 9).Value = 1605
$ws.Cells.Item(16, 11).Value = 1605
$ws.Cells.Item(16, 13).Value = -1318
$ws.Cells.Item(22, 8).Value = 199.04
$ws.Cells.Item(22, 9).Value = 183.22449
$ws.Cells.Item(22, 10).Value = 974
$ws.Cells.Item(22, 11).Value = 183.22449
$ws.Cells.Item(22, 12).Value = 974
$ws.Cells.Item(22, 13).Value = 166.77551
$ws.Cells.Item(22, 14).Value = -1674
$ws.Cells.Item(38, 8).Value = 10000
$ws.Cells.Item(38, 10).Value = 10000
$ws.Cells.Item(38, 12).Value = 10000
$ws.Cells.Item(38, 14).Value = -10754
$ws.Cells.Item(42, 8).Value = 12099.667
$ws.Cells.Item(42, 10).Value = 18000
$ws.Cells.Item(42, 12).Value = 18000
$ws.Cells.Item(42, 14).Value = -19186
$ws.Cells.Item(46, 8).Value = 10000
$ws.Cells.Item(46, 10).Value = 10000
$ws.Cells.Item(46, 12).Value = 10000
$ws.Cells.Item(46, 14).Value = -10422
$ws.Cells.Item(55, 8).Value = 14449
$ws.Cells.Item(55, 9).Value = 14449
$ws.Cells.Item(55, 11).Value = 14449
$ws.Cells.Item(55, 13).Value = -14134
$ws.Cells.Item(56, 8).Value = 42000
$ws.Cells.Item(56, 10).Value = 59000
$ws.Cells.Item(56, 12).Value = 59000
$ws.Cells.Item(56, 14).Value = -60690
$ws.Cells.Item(58, 8).Value = 3246.8333
$ws.Cells.Item(58, 9).Value = 1827
$ws.Cells.Item(58, 10).Value = 4666.6665
$ws.Cells.Item(58, 11).Value = 1827
$ws.Cells.Item(58, 12).Value = 4666.6665
$ws.Cells.Item(58, 13).Value = -1624
$ws.Cells.Item(58, 14).Value = -5072.6665
$ws.Cells.Item(82, 8).Value = 56000
$ws.Cells.Item(82, 9).Value = 56000
$ws.Cells.Item(82, 11).Value = 56000
$ws.Cells.Item(82, 13).Value = -55639
$ws.Cells.Item(85, 8).Value = 56000
$ws.Cells.Item(85, 9).Value = 56000
$ws.Cells.Item(85, 11).Value = 56000
$ws.Cells.Item(85, 13).Value = -54752
$ws.Cells.Item(105, 8).Value = 949.5
$ws.Cells.Item(105, 9).Value = 949.5
$ws.Cells.Item(105, 11).Value = 949.5
$ws.Cells.Item(105, 13).Value = 797.5
$ws.Cells.Item(113, 8).Value = 1703.3334
$ws.Cells.Item(113, 9).Value = 1605
$ws.Cells.Item(113, 11).Value = 1605
$ws.Cells.Item(113, 13).Value = 565
$ws.Cells.Item(122, 8).Value = 9500
$ws.Cells.Item(122, 9).Value = 9500
$ws.Cells.Item(122, 11).Value = 28500
$ws.Cells.Item(122, 13).Value = -26050
$ws.Cells.Item(134, 8).Value = 9526747
$ws.Cells.Item(134, 9).Value = 12989201
$ws.Cells.Item(134, 10).Value = 4999.75
$ws.Cells.Item(134, 11).Value = 38967603
$ws.Cells.Item(134, 12).Value = 14999.25
$ws.Cells.Item(134, 13).Value = -38965068
$ws.Cells.Item(134, 14).Value = -20069.25
$ws.Cells.Item(136, 8).Value = 3246.8333
$ws.Cells.Item(136, 9).Value = 1827
$ws.Cells.Item(136, 10).Value = 4666.6665
$ws.Cells.Item(136, 11).Value = 5481
$ws.Cells.Item(136, 12).Value = 13999.9995
$ws.Cells.Item(136, 13).Value = -2931
$ws.Cells.Item(136, 14).Value = -19099.9995
$ws.Cells.Item(140, 8).Value = 42500
$ws.Cells.Item(140, 10).Value = 42500
$ws.Cells.Item(140, 12).Value = 42500
$ws.Cells.Item(140, 14).Value = -52860

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(9, 8).Value = 7143
$ws.Cells.Item(9, 10).Value = 9600.200000000001
$ws.Cells.Item(9, 12).Value = 28800.6
$ws.Cells.Item(9, 14).Value = -29248.6
$ws.Cells.Item(34, 8).Value = 4822
$ws.Cells.Item(34, 10).Value = 5343.625
$ws.Cells.Item(34, 12).Value = 16030.875
$ws.Cells.Item(34, 14).Value = -16198.875
$ws.Cells.Item(39, 8).Value = 8400
$ws.Cells.Item(39, 10).Value = 8400
$ws.Cells.Item(39, 12).Value = 25200
$ws.Cells.Item(39, 14).Value = -25788
$ws.Cells.Item(55, 8).Value = 557
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 14).Value = $null
$ws.Cells.Item(107, 8).Value = 356.75
$ws.Cells.Item(107, 9).Value = 272
$ws.Cells.Item(107, 10).Value = 441.5
$ws.Cells.Item(107, 11).Value = 816
$ws.Cells.Item(107, 12).Value = 1324.5
$ws.Cells.Item(107, 13).Value = 1104
$ws.Cells.Item(107, 14).Value = -5164.5
$ws.Cells.Item(133, 8).Value = 7000
$ws.Cells.Item(133, 9).Value = 7000
$ws.Cells.Item(133, 11).Value = 21000
$ws.Cells.Item(133, 13).Value = -15940

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20, 8).Value = 15000
$ws.Cells.Item(20, 10).Value = 15000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 14).Value = -15490
$ws.Cells.Item(24, 8).Value = 22824.928
$ws.Cells.Item(24, 9).Value = 19195.77
$ws.Cells.Item(24, 11).Value = 19195.77
$ws.Cells.Item(24, 13).Value = -19022.77
$ws.Cells.Item(102, 8).Value = 2316
$ws.Cells.Item(102, 9).Value = 2316
$ws.Cells.Item(102, 11).Value = 2316
$ws.Cells.Item(102, 13).Value = -694
$ws.Cells.Item(122, 8).Value = 2482.3076
$ws.Cells.Item(122, 9).Value = 2482.3076
$ws.Cells.Item(122, 11).Value = 7446.9228
$ws.Cells.Item(122, 13).Value = -4996.9228
$ws.Cells.Item(132, 8).Value = 3443.4285
$ws.Cells.Item(132, 9).Value = 3567.3333
$ws.Cells.Item(132, 10).Value = 2700
$ws.Cells.Item(132, 11).Value = 10701.9999
$ws.Cells.Item(132, 12).Value = 8100
$ws.Cells.Item(132, 13).Value = -8171.999899999999
$ws.Cells.Item(132, 14).Value = -13160
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).Value = $null
$ws.Cells.Item(134, 8).Value = 49495.5
$ws.Cells.Item(134, 10).Value = 49495.5
$ws.Cells.Item(134, 12).Value = 148486.5
$ws.Cells.Item(134, 14).Value = -153556.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 99997
$ws.Cells.Item(36, 10).Value = 99997
$ws.Cells.Item(36, 12).Value = 99997
$ws.Cells.Item(36, 14).Value = -101121
$ws.Cells.Item(40, 8).Value = 17825
$ws.Cells.Item(40, 9).Value = 15704.091
$ws.Cells.Item(40, 11).Value = 15704.091
$ws.Cells.Item(40, 13).Value = -15568.091
$ws.Cells.Item(93, 8).Value = 842.5714
$ws.Cells.Item(93, 9).Value = 733
$ws.Cells.Item(93, 10).Value = 924.75
$ws.Cells.Item(93, 11).Value = 733
$ws.Cells.Item(93, 12).Value = 924.75
$ws.Cells.Item(93, 13).Value = 515
$ws.Cells.Item(93, 14).Value = -3420.75
$ws.Cells.Item(132, 8).Value = 2646.1428
$ws.Cells.Item(132, 9).Value = 2551.3333
$ws.Cells.Item(132, 10).Value = 2717.25
$ws.Cells.Item(132, 11).Value = 7653.999899999999
$ws.Cells.Item(132, 12).Value = 8151.75
$ws.Cells.Item(132, 13).Value = -5123.999899999999
$ws.Cells.Item(132, 14).Value = -13211.75
$ws.Cells.Item(136, 8).Value = 1410.75
$ws.Cells.Item(136, 9).Value = 1247.6666
$ws.Cells.Item(136, 11).Value = 3742.9998
$ws.Cells.Item(136, 13).Value = -1192.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 43333
$ws.Cells.Item(45, 9).Value = 50000
$ws.Cells.Item(45, 11).Value = 50000
$ws.Cells.Item(45, 13).Value = -49509
$ws.Cells.Item(107, 8).Value = 3871
$ws.Cells.Item(107, 9).Value = 3762.077
$ws.Cells.Item(107, 10).Value = 4225
$ws.Cells.Item(107, 11).Value = 11286.231
$ws.Cells.Item(107, 12).Value = 12675
$ws.Cells.Item(107, 13).Value = -9366.231
$ws.Cells.Item(107, 14).Value = -16515
$ws.Cells.Item(132, 8).Value = 4847.364
$ws.Cells.Item(132, 9).Value = 3826.5
$ws.Cells.Item(132, 10).Value = 7569.6665
$ws.Cells.Item(132, 11).Value = 11479.5
$ws.Cells.Item(132, 12).Value = 22708.9995
$ws.Cells.Item(132, 13).Value = -8949.5
$ws.Cells.Item(132, 14).Value = -27768.9995
$ws.Cells.Item(136, 8).Value = 7222.3335
$ws.Cells.Item(136, 9).Value = 6302.5
$ws.Cells.Item(136, 11).Value = 18907.5
$ws.Cells.Item(136, 13).Value = -16357.5
